# Split the "Programa" section's two run-of-text paragraphs (PT and EN)
# into several sentences separated by manual line breaks (<w:br/>),
# by using Find/Replace with the "^l" (manual line break) special
# sequence at each semicolon boundary between topics.

$d = $word.ActiveDocument

# --- Portuguese paragraph ---
$d.Content.Find.Execute(
    "(mínimo 3); Inovação", $true, $false, $false, $false, $false, $true, 1, $false,
    "(mínimo 3);^l Inovação", 2) | Out-Null

$d.Content.Find.Execute(
    "características;Legislação", $true, $false, $false, $false, $false, $true, 1, $false,
    "características;^lLegislação", 2) | Out-Null

$d.Content.Find.Execute(
    "empresarial;Gerenciamento", $true, $false, $false, $false, $false, $true, 1, $false,
    "empresarial;^lGerenciamento", 2) | Out-Null

$d.Content.Find.Execute(
    "causas;Formulação", $true, $false, $false, $false, $false, $true, 1, $false,
    "causas;^lFormulação", 2) | Out-Null

$d.Content.Find.Execute(
    "etc;Especificação", $true, $false, $false, $false, $false, $true, 1, $false,
    "etc;^lEspecificação", 2) | Out-Null

$d.Content.Find.Execute(
    "Decisão;Elaboração", $true, $false, $false, $false, $false, $true, 1, $false,
    "Decisão;^lElaboração", 2) | Out-Null

# --- English (italic) paragraph ---
$d.Content.Find.Execute(
    "(minimum 3);Systematic", $true, $false, $false, $false, $false, $true, 1, $false,
    "(minimum 3);^lSystematic", 2) | Out-Null

$d.Content.Find.Execute(
    "characteristics;Legislation", $true, $false, $false, $false, $false, $true, 1, $false,
    "characteristics;^lLegislation", 2) | Out-Null

$d.Content.Find.Execute(
    "action;Project and Schedule", $true, $false, $false, $false, $false, $true, 1, $false,
    "action;^lProject and Schedule", 2) | Out-Null

$d.Content.Find.Execute(
    "causes;Project Formulation", $true, $false, $false, $false, $false, $true, 1, $false,
    "causes;^lProject Formulation", 2) | Out-Null

$d.Content.Find.Execute(
    "etc;Problem Specification", $true, $false, $false, $false, $false, $true, 1, $false,
    "etc;^lProblem Specification", 2) | Out-Null

$d.Content.Find.Execute(
    "making;Preparation of reports", $true, $false, $false, $false, $false, $true, 1, $false,
    "making;^lPreparation of reports", 2) | Out-Null
